# Auto-generated edit script applying per-cell value updates described by the diff.
# Workbook: FFXIV Leve profit tracker ("Phantom_Profits"), 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# each backed by a Table (Table_<sheet>) with columns A-N. Columns H-N (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) are refreshed market-price-derived values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2033
$ws.Range("I62").Value = 2033
$ws.Range("K62").Value = 2033
$ws.Range("M62").Value = -1409

$ws.Range("H65").Value = 2033
$ws.Range("I65").Value = 2033
$ws.Range("K65").Value = 10165
$ws.Range("M65").Value = -7045

$ws.Range("H112").Value = 3018.2942
$ws.Range("J112").Value = 3082.875
$ws.Range("L112").Value = 9248.625
$ws.Range("N112").Value = -11464.625

$ws.Range("H116").Value = 4501.143
$ws.Range("I116").Value = 3126.5
$ws.Range("J116").Value = 6334
$ws.Range("K116").Value = 3126.5
$ws.Range("L116").Value = 6334
$ws.Range("M116").Value = 315.5
$ws.Range("N116").Value = -13218

$ws.Range("H137").Value = 22224462
$ws.Range("I137").Value = 66667580
$ws.Range("J137").Value = 2902.7
$ws.Range("K137").Value = 200002740
$ws.Range("L137").Value = 8708.099999999999
$ws.Range("M137").Value = -200000190
$ws.Range("N137").Value = -13808.1

$ws.Range("H138").Value = 8264.177
$ws.Range("J138").Value = 9902.056
$ws.Range("L138").Value = 29706.168
$ws.Range("N138").Value = -39986.16800000001

$ws.Range("H139").Value = 69997.5
$ws.Range("J139").Value = 69997.5
$ws.Range("L139").Value = 69997.5
$ws.Range("N139").Value = -80277.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5306.64
$ws.Range("I61").Value = 4855.087
$ws.Range("J61").Value = 10499.5
$ws.Range("K61").Value = 4855.087
$ws.Range("L61").Value = 10499.5
$ws.Range("M61").Value = -4643.087
$ws.Range("N61").Value = -10923.5

$ws.Range("H92").Value = 152966.33
$ws.Range("J92").Value = 152966.33
$ws.Range("L92").Value = 152966.33
$ws.Range("N92").Value = -157958.33

$ws.Range("H136").Value = 5306.64
$ws.Range("I136").Value = 4855.087
$ws.Range("J136").Value = 10499.5
$ws.Range("K136").Value = 14565.261
$ws.Range("L136").Value = 31498.5
$ws.Range("M136").Value = -12015.261
$ws.Range("N136").Value = -36598.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1625.2858
$ws.Range("I105").Value = 1523.3636
$ws.Range("J105").Value = 1999
$ws.Range("K105").Value = 1523.3636
$ws.Range("L105").Value = 1999
$ws.Range("M105").Value = 223.6364000000001
$ws.Range("N105").Value = -5493

$ws.Range("H107").Value = 1478.4286
$ws.Range("I107").Value = 1478.4286
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1478.4286
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 441.5714
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 3500
$ws.Range("I23").Value = 3500
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 3500
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -3260
$ws.Range("N23").ClearContents()

$ws.Range("H27").Value = 3500
$ws.Range("I27").Value = 3500
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 3500
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -3308
$ws.Range("N27").ClearContents()

$ws.Range("H31").Value = 2734.7222
$ws.Range("I31").Value = 2622.0667
$ws.Range("J31").Value = 3298
$ws.Range("K31").Value = 2622.0667
$ws.Range("L31").Value = 3298
$ws.Range("M31").Value = -2327.0667
$ws.Range("N31").Value = -3888

$ws.Range("H34").Value = 2734.7222
$ws.Range("I34").Value = 2622.0667
$ws.Range("J34").Value = 3298
$ws.Range("K34").Value = 2622.0667
$ws.Range("L34").Value = 3298
$ws.Range("M34").Value = -2420.0667
$ws.Range("N34").Value = -3702

$ws.Range("H86").Value = 6240.6665
$ws.Range("I86").Value = 6240.6665
$ws.Range("K86").Value = 6240.6665
$ws.Range("M86").Value = -5117.6665

$ws.Range("H89").Value = 6240.6665
$ws.Range("I89").Value = 6240.6665
$ws.Range("K89").Value = 31203.3325
$ws.Range("M89").Value = -25587.3325

$ws.Range("H107").Value = 762.4375
$ws.Range("I107").Value = 558.3333
$ws.Range("J107").Value = 1374.75
$ws.Range("K107").Value = 558.3333
$ws.Range("L107").Value = 1374.75
$ws.Range("M107").Value = 1361.6667
$ws.Range("N107").Value = -5214.75

$ws.Range("H132").Value = 22237024
$ws.Range("I132").Value = 28589544
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 85768632
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -85766102
$ws.Range("N132").Value = -14660

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 2000
$ws.Range("I116").Value = 2000
$ws.Range("K116").Value = 6000
$ws.Range("M116").Value = -2558

$ws.Range("H137").Value = 55000
$ws.Range("J137").Value = 10000
$ws.Range("L137").Value = 30000
$ws.Range("N137").Value = -40200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6047.625
$ws.Range("I70").Value = 6231.1665
$ws.Range("K70").Value = 6231.1665
$ws.Range("M70").Value = -5961.1665

$ws.Range("H73").Value = 6047.625
$ws.Range("I73").Value = 6231.1665
$ws.Range("K73").Value = 6231.1665
$ws.Range("M73").Value = -5295.1665

$ws.Range("H132").Value = 111115016
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 111115016
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 333345048
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -333350108

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2468.3076
$ws.Range("I132").Value = 2465.6667
$ws.Range("K132").Value = 7397.000100000001
$ws.Range("M132").Value = -4867.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9900
$ws.Range("I62").Value = 9900
$ws.Range("K62").Value = 9900
$ws.Range("M62").Value = -9276

$ws.Range("H65").Value = 9900
$ws.Range("I65").Value = 9900
$ws.Range("K65").Value = 49500
$ws.Range("M65").Value = -46380

$ws.Range("H69").Value = 28141.3
$ws.Range("J69").Value = 26240.777
$ws.Range("L69").Value = 26240.777
$ws.Range("N69").Value = -27738.777

$ws.Range("H72").Value = 28141.3
$ws.Range("J72").Value = 26240.777
$ws.Range("L72").Value = 78722.33099999999
$ws.Range("N72").Value = -86210.33099999999

$ws.Range("H122").Value = 5605.25
$ws.Range("I122").Value = 5605.25
$ws.Range("K122").Value = 16815.75
$ws.Range("M122").Value = -14365.75

$ws.Range("H132").Value = 200004000
$ws.Range("I132").Value = 4999.25
$ws.Range("K132").Value = 14997.75
$ws.Range("M132").Value = -12467.75

$ws.Range("H136").Value = 9029.299999999999
$ws.Range("I136").Value = 10989.692
$ws.Range("K136").Value = 32969.076
$ws.Range("M136").Value = -30419.076
